$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, border, centered) from E1 into F1, then set header text
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# Fill time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 10:51:13.155319"
$ws.Range("F3").Value = "2021-10-05 10:51:13.155331"
$ws.Range("F4").Value = "2021-10-05 10:51:13.155334"
$ws.Range("F5").Value = "2021-10-05 10:51:13.155337"
$ws.Range("F6").Value = "2021-10-05 10:51:13.155340"
$ws.Range("F7").Value = "2021-10-05 10:51:13.155342"
$ws.Range("F8").Value = "2021-10-05 10:51:13.155345"
$ws.Range("F9").Value = "2021-10-05 10:51:13.155347"
$ws.Range("F10").Value = "2021-10-05 10:51:13.155350"
$ws.Range("F11").Value = "2021-10-05 10:51:13.155353"
$ws.Range("F12").Value = "2021-10-05 10:51:13.155355"
$ws.Range("F13").Value = "2021-10-05 10:51:13.155358"
$ws.Range("F14").Value = "2021-10-05 10:51:13.155360"
$ws.Range("F15").Value = "2021-10-05 10:51:13.155363"
$ws.Range("F16").Value = "2021-10-05 10:51:13.155365"
$ws.Range("F17").Value = "2021-10-05 10:51:13.155368"
$ws.Range("F18").Value = "2021-10-05 10:51:13.155370"
$ws.Range("F19").Value = "2021-10-05 10:51:13.155373"
$ws.Range("F20").Value = "2021-10-05 10:51:13.155376"
$ws.Range("F21").Value = "2021-10-05 10:51:13.155378"
$ws.Range("F22").Value = "2021-10-05 10:51:13.155381"
$ws.Range("F23").Value = "2021-10-05 10:51:13.155383"
$ws.Range("F24").Value = "2021-10-05 10:51:13.155385"
$ws.Range("F25").Value = "2021-10-05 10:51:13.155388"
$ws.Range("F26").Value = "2021-10-05 10:51:13.155391"
$ws.Range("F27").Value = "2021-10-05 10:51:13.155393"
$ws.Range("F28").Value = "2021-10-05 10:51:13.155396"
$ws.Range("F29").Value = "2021-10-05 10:51:13.155398"
$ws.Range("F30").Value = "2021-10-05 10:51:13.155401"
$ws.Range("F31").Value = "2021-10-05 10:51:13.155403"
$ws.Range("F32").Value = "2021-10-05 10:51:13.155406"
$ws.Range("F33").Value = "2021-10-05 10:51:13.155408"
$ws.Range("F34").Value = "2021-10-05 10:51:13.155411"
$ws.Range("F35").Value = "2021-10-05 10:51:13.155413"
$ws.Range("F36").Value = "2021-10-05 10:51:13.155416"
$ws.Range("F37").Value = "2021-10-05 10:51:13.155418"
$ws.Range("F38").Value = "2021-10-05 10:51:13.155421"
$ws.Range("F39").Value = "2021-10-05 10:51:13.155424"
$ws.Range("F40").Value = "2021-10-05 10:51:13.155426"
$ws.Range("F41").Value = "2021-10-05 10:51:13.155429"
$ws.Range("F42").Value = "2021-10-05 10:51:13.155432"
$ws.Range("F43").Value = "2021-10-05 10:51:13.155435"
$ws.Range("F44").Value = "2021-10-05 10:51:13.155437"
$ws.Range("F45").Value = "2021-10-05 10:51:13.155440"
$ws.Range("F46").Value = "2021-10-05 10:51:13.155442"
$ws.Range("F47").Value = "2021-10-05 10:51:13.155445"
$ws.Range("F48").Value = "2021-10-05 10:51:13.155447"
$ws.Range("F49").Value = "2021-10-05 10:51:13.155449"
$ws.Range("F50").Value = "2021-10-05 10:51:13.155452"
$ws.Range("F51").Value = "2021-10-05 10:51:13.155454"
$ws.Range("F52").Value = "2021-10-05 10:51:13.155457"
$ws.Range("F53").Value = "2021-10-05 10:51:13.155459"
$ws.Range("F54").Value = "2021-10-05 10:51:13.155462"

Write-Output "Done"
